# Update results for Steel (Sheet1 of 2030_FR.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hydrogen row (B3): updated value
$ws.Range("B3").Value = 1354.070613502361

# Methanol row (C4): updated value
$ws.Range("C4").Value = 170.2654775220201

# Other row (D8): updated value
$ws.Range("D8").Value = 1126.966996040081

$wb.Save()
